$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 24461.5359505928
$ws.Range("E2").Value = 38965.3981072151
$ws.Range("F2").Value = 50045.3189749633
$ws.Range("I2").Value = 7302.53595059279
$ws.Range("B3").Value = 13751.7745901086
$ws.Range("E3").Value = 18937.4695818333
$ws.Range("F3").Value = 22084.6511996782
$ws.Range("I3").Value = 5733.77459010856
$ws.Range("B4").Value = 18538.4450769877
$ws.Range("C4").Value = 12694.6554253341
$ws.Range("D4").Value = 11441.6264016489
$ws.Range("I4").Value = 12510.4450769877
$ws.Range("B5").Value = 27188.5952362007
$ws.Range("C5").Value = 14987.1804181161
$ws.Range("D5").Value = 12780.4867374626
$ws.Range("I5").Value = 21019.5952362007
$ws.Range("B6").Value = 36446.5684537737
$ws.Range("C6").Value = 17277.2105196283
$ws.Range("D6").Value = 14647.7331008109
$ws.Range("I6").Value = 27173.5684537737
$ws.Range("B7").Value = 39424.2224833897
$ws.Range("C7").Value = 19205.41845331
$ws.Range("D7").Value = 16018.2984224765
$ws.Range("I7").Value = 27429.2224833897
$ws.Range("B8").Value = 33060.5906786407
$ws.Range("I8").Value = 24052.5906786407
$ws.Range("B9").Value = 21377.2861751258
$ws.Range("E9").Value = 36693.2319133197
$ws.Range("I9").Value = 12369.2861751258
$ws.Range("B10").Value = 19533.3406833496
$ws.Range("E10").Value = 34850.3596412676
$ws.Range("I10").Value = 6321.3406833496
$ws.Range("B11").Value = 20556.5407096913
$ws.Range("E11").Value = 36765.7237617603
$ws.Range("I11").Value = 7738.5407096913
$ws.Range("B12").Value = 25092.9864667262
$ws.Range("I12").Value = 9881.9864667262
$ws.Range("B13").Value = 27728.3048946825
$ws.Range("I13").Value = 16291.3048946825
$ws.Range("B14").Value = 24462.2257899274
$ws.Range("I14").Value = 17384.2257899274
$ws.Range("B15").Value = 13993.2920688145
$ws.Range("E15").Value = 28524.5981673961
$ws.Range("F15").Value = 35975.9047609483
$ws.Range("I15").Value = 9600.29206881454
$ws.Range("B16").Value = 18464.1487373481
$ws.Range("I16").Value = 8860.14873734811
$ws.Range("B17").Value = 26589.0912794963
$ws.Range("C17").Value = 13277.8648708032
$ws.Range("I17").Value = 15088.0912794963
$ws.Range("B18").Value = 35140.122205063
$ws.Range("C18").Value = 14944.719734274
$ws.Range("D18").Value = 12425.8685165999
$ws.Range("I18").Value = 23170.122205063
$ws.Range("B19").Value = 37768.5659482412
$ws.Range("C19").Value = 16211.8616118503
$ws.Range("I19").Value = 25753.5659482412
$ws.Range("B20").Value = 31906.7614493348
$ws.Range("I20").Value = 22539.7614493348
$ws.Range("B21").Value = 21499.1943208056
$ws.Range("I21").Value = 13146.1943208056
$ws.Range("B22").Value = 20100.8182111201
$ws.Range("E22").Value = 43268.4953080088
$ws.Range("I22").Value = 7494.8182111201
$ws.Range("B23").Value = 21104.3450208666
$ws.Range("I23").Value = 9988.34502086657
$ws.Range("B24").Value = 25260.6795519576
$ws.Range("I24").Value = 13379.6795519576
$ws.Range("B25").Value = 27736.6424954553
$ws.Range("I25").Value = 16844.6424954553
$ws.Range("B26").Value = 24496.9390469305
$ws.Range("I26").Value = 17926.9390469305
$ws.Range("B27").Value = 15524.3032770952
$ws.Range("E27").Value = 38751.2723385486
$ws.Range("F27").Value = 53814.0992316784
$ws.Range("I27").Value = 11033.3032770952
$ws.Range("B28").Value = 19333.023764692
$ws.Range("I28").Value = 9332.02376469202
$ws.Range("B29").Value = 26950.9089810644
$ws.Range("I29").Value = 16982.9089810644
$ws.Range("B30").Value = 35088.2275010212
$ws.Range("C30").Value = 13576.3747642065
$ws.Range("I30").Value = 23937.2275010212
$ws.Range("B31").Value = 37415.2258399519
$ws.Range("C31").Value = 14504.5921919688
$ws.Range("I31").Value = 25470.2258399519
$ws.Range("B32").Value = 31418.0556653882
$ws.Range("I32").Value = 22027.0556653882
$ws.Range("B33").Value = 21255.9454386756
$ws.Range("I33").Value = 12783.9454386756
$ws.Range("B34").Value = 20062.3624705834
$ws.Range("I34").Value = 9021.36247058342
$ws.Range("B35").Value = 21183.5275991687
$ws.Range("I35").Value = 11646.5275991687
$ws.Range("B36").Value = 25365.3760418774
$ws.Range("I36").Value = 16663.3760418774
$ws.Range("B37").Value = 27822.4130179549
$ws.Range("I37").Value = 23983.4130179549
$ws.Range("B38").Value = 24597.4011873646
$ws.Range("I38").Value = 22227.4011873646
$ws.Range("B39").Value = 16839.4729831057
$ws.Range("E39").Value = 44451.2432879007
$ws.Range("F39").Value = 68826.5644353505
$ws.Range("I39").Value = 12291.4729831057
$ws.Range("B40").Value = 20532.3198370889
$ws.Range("I40").Value = 13233.3198370889
$ws.Range("B41").Value = 27867.1707166316
$ws.Range("I41").Value = 20839.1707166316
$ws.Range("B42").Value = 35701.3579121552
$ws.Range("C42").Value = 13048.548602695
$ws.Range("I42").Value = 26771.3579121552
$ws.Range("B43").Value = 37749.011355913
$ws.Range("C43").Value = 14045.4937767376
$ws.Range("I43").Value = 27039.011355913
$ws.Range("B44").Value = 31639.2901838467
$ws.Range("I44").Value = 22359.2901838467
$ws.Range("B45").Value = 21378.5771325739
$ws.Range("I45").Value = 13459.5771325739
$ws.Range("B46").Value = 20143.478384449
$ws.Range("I46").Value = 9276.47838444903
$ws.Range("B47").Value = 21269.8883947069
$ws.Range("I47").Value = 12982.8883947069
$ws.Range("B48").Value = 25458.6812323659
$ws.Range("I48").Value = 17816.6812323659
$ws.Range("B49").Value = 27919.9019948061
$ws.Range("I49").Value = 20827.9019948061
